$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 0.2182062694869273
$ws.Cells.Item(2, 3).Value = 0.727114017474318
$ws.Cells.Item(2, 4).Value = 0.3056242433470812
$ws.Cells.Item(3, 2).Value = 0.1728940631576693
$ws.Cells.Item(3, 3).Value = 0.3189668226482641
$ws.Cells.Item(3, 4).Value = 0.2806941821307254
$ws.Cells.Item(4, 2).Value = 0.2026979580823564
$ws.Cells.Item(4, 3).Value = 1.988221847353749
$ws.Cells.Item(4, 4).Value = 0.3139868914312781
$ws.Cells.Item(5, 2).Value = 0.2135013263428866
$ws.Cells.Item(5, 3).Value = 1.466772960679729
$ws.Cells.Item(5, 4).Value = 0.2111532367743074
$ws.Cells.Item(6, 2).Value = 0.1790377367902297
$ws.Cells.Item(6, 3).Value = 3.680984438311768
$ws.Cells.Item(6, 4).Value = 0.2544067902771067
$ws.Cells.Item(7, 2).Value = 0.1476591858533107
$ws.Cells.Item(7, 3).Value = 4.032918071143688
$ws.Cells.Item(7, 4).Value = 0.2311548327076659
$ws.Cells.Item(8, 2).Value = 0.1689252967058519
$ws.Cells.Item(8, 3).Value = 6.002255972224708
$ws.Cells.Item(8, 4).Value = 0.2627284189232876
$ws.Cells.Item(9, 2).Value = 0.2253248938005501
$ws.Cells.Item(9, 3).Value = 8.616201250615793
$ws.Cells.Item(9, 4).Value = 0.183676005822393
$ws.Cells.Item(10, 2).Value = 0.2840942207071662
$ws.Cells.Item(10, 3).Value = 10.77316973899689
$ws.Cells.Item(10, 4).Value = 0.4144955855970494
$ws.Cells.Item(11, 2).Value = 0.1157406564722776
$ws.Cells.Item(11, 3).Value = 14.39979474203238
$ws.Cells.Item(11, 4).Value = 0.2917191737295326
$ws.Cells.Item(12, 2).Value = 0.1308992366970819
$ws.Cells.Item(12, 3).Value = 16.79950491845142
$ws.Cells.Item(12, 4).Value = 0.3647787250130659
$ws.Cells.Item(13, 2).Value = 0.1188024213760002
$ws.Cells.Item(13, 3).Value = 20.50627134335875
$ws.Cells.Item(13, 4).Value = 0.28004431159937
$ws.Cells.Item(14, 2).Value = 0.1772004303011717
$ws.Cells.Item(14, 3).Value = 24.13902237633232
$ws.Cells.Item(14, 4).Value = 0.2559651345863
$ws.Cells.Item(15, 2).Value = 0.1358568744455128
$ws.Cells.Item(15, 3).Value = 29.6106011778766
$ws.Cells.Item(15, 4).Value = 0.3122148700369549
$ws.Cells.Item(16, 2).Value = 0.1585866003157774
$ws.Cells.Item(16, 3).Value = 33.76421156042156
$ws.Cells.Item(16, 4).Value = 0.21988898692052
$ws.Cells.Item(17, 2).Value = 0.1544096646683313
$ws.Cells.Item(17, 3).Value = 38.70524170350165
$ws.Cells.Item(17, 4).Value = 0.2559563306151532
$ws.Cells.Item(18, 2).Value = 0.2382661534364694
$ws.Cells.Item(18, 3).Value = 44.61726158553743
$ws.Cells.Item(18, 4).Value = 0.3879057784949178
$ws.Cells.Item(19, 2).Value = 0.1398407872048483
$ws.Cells.Item(19, 3).Value = 50.39233091677868
$ws.Cells.Item(19, 4).Value = 0.2216534357382231
